# ABA Customer Info.xlsx — "all pages are finished"
# Adds Country / Gender / Email (hyperlinked) / Password columns, replaces the
# old "Reservation Number" column with "Country", updates Drivers License #
# and Age for several rows, and re-points the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Drivers License # (D) updates
# ---------------------------------------------------------------------------
$ws.Cells.Item(2,4).Value  = 123456789
$ws.Cells.Item(3,4).Value  = 133456789
$ws.Cells.Item(4,4).Value  = 986754382
$ws.Cells.Item(5,4).Value  = 364527183
$ws.Cells.Item(6,4).Value  = 321321312
$ws.Cells.Item(7,4).Value  = 654654654
$ws.Cells.Item(8,4).Value  = 888999000
$ws.Cells.Item(9,4).Value  = 970960593
$ws.Cells.Item(10,4).Value = 563018920

# ---------------------------------------------------------------------------
# 2. Age (E) updates — only the rows that actually changed
# ---------------------------------------------------------------------------
$ws.Cells.Item(4,5).Value = 21
$ws.Cells.Item(6,5).Value = 19
$ws.Cells.Item(8,5).Value = 22

# ---------------------------------------------------------------------------
# 3. Column F: "Reservation Number" -> "Country"
# ---------------------------------------------------------------------------
$ws.Cells.Item(1,6).Value  = "Country"
$ws.Cells.Item(3,6).Value  = "China"
$ws.Cells.Item(4,6).Value  = "Canada"
$ws.Cells.Item(2,6).Value  = "England"
$ws.Cells.Item(6,6).Value  = "Japan"
$ws.Cells.Item(8,6).Value  = "United States"
$ws.Cells.Item(9,6).Value  = "Mexico"
$ws.Cells.Item(10,6).Value = "Hong Kong"
$ws.Cells.Item(5,6).Value  = "Canada"
$ws.Cells.Item(7,6).Value  = "Canada"

# ---------------------------------------------------------------------------
# 4. Column G: Gender
# ---------------------------------------------------------------------------
$ws.Cells.Item(1,7).Value  = "Gender"
$ws.Cells.Item(2,7).Value  = "F"
$ws.Cells.Item(3,7).Value  = "M"
$ws.Cells.Item(4,7).Value  = "M"
$ws.Cells.Item(5,7).Value  = "M"
$ws.Cells.Item(6,7).Value  = "F"
$ws.Cells.Item(7,7).Value  = "F"
$ws.Cells.Item(8,7).Value  = "F"
$ws.Cells.Item(9,7).Value  = "F"
$ws.Cells.Item(10,7).Value = "M"

# ---------------------------------------------------------------------------
# 5. Column I header: Password (typed before the Email header)
# ---------------------------------------------------------------------------
$ws.Cells.Item(1,9).Value = "Password"

# ---------------------------------------------------------------------------
# 6. Column H: Email (header + values). Most rows become real mailto:
#    hyperlinks; rows 3 and 10 stay plain text (no hyperlink / no style).
# ---------------------------------------------------------------------------
$ws.Cells.Item(1,8).Value = "Email"

$ws.Cells.Item(2,8).Value = "alana@nomail.com"
$ws.Hyperlinks.Add($ws.Cells.Item(2,8), "mailto:alana@nomail.com") | Out-Null

$ws.Cells.Item(10,8).Value = "archie@nomail.com"

$ws.Cells.Item(9,8).Value = "betty@nomail.com"
$ws.Hyperlinks.Add($ws.Cells.Item(9,8), "mailto:betty@nomail.com") | Out-Null

$ws.Cells.Item(8,8).Value = "veronica@nomail.com"
$ws.Hyperlinks.Add($ws.Cells.Item(8,8), "mailto:veronica@nomail.com") | Out-Null

$ws.Cells.Item(7,8).Value = "brianna@nomail.com"
$ws.Hyperlinks.Add($ws.Cells.Item(7,8), "mailto:brianna@nomail.com") | Out-Null

$ws.Cells.Item(6,8).Value = "chrissy@nomail.com"
$ws.Hyperlinks.Add($ws.Cells.Item(6,8), "mailto:chrissy@nomail.com") | Out-Null

$ws.Cells.Item(5,8).Value = "derrick@nomail.com"
$ws.Hyperlinks.Add($ws.Cells.Item(5,8), "mailto:derrick@nomail.com") | Out-Null

$ws.Cells.Item(4,8).Value = "cooper@nomail.com"
$ws.Hyperlinks.Add($ws.Cells.Item(4,8), "mailto:cooper@nomail.com") | Out-Null

$ws.Cells.Item(3,8).Value = "alan@nomail.com"

# ---------------------------------------------------------------------------
# 7. Column I: Password values
# ---------------------------------------------------------------------------
$ws.Cells.Item(2,9).Value  = 3214
$ws.Cells.Item(3,9).Value  = 3232
$ws.Cells.Item(4,9).Value  = 5454
$ws.Cells.Item(5,9).Value  = 6565
$ws.Cells.Item(6,9).Value  = 7676
$ws.Cells.Item(7,9).Value  = 4344
$ws.Cells.Item(8,9).Value  = 9796
$ws.Cells.Item(9,9).Value  = 6043
$ws.Cells.Item(10,9).Value = 2389

# ---------------------------------------------------------------------------
# 8. Column widths (new H/I columns + minor width touch-ups on existing ones)
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 21.666666666666668
$ws.Columns.Item(2).ColumnWidth = 32.330729166666664
$ws.Columns.Item(4).ColumnWidth = 33.998697916666664
$ws.Columns.Item(6).ColumnWidth = 36.998697916666664
$ws.Columns.Item(8).ColumnWidth = 17.498697916666668
$ws.Columns.Item(9).ColumnWidth = 12.498697916666666

# ---------------------------------------------------------------------------
# 9. Selection / active cell
# ---------------------------------------------------------------------------
$ws.Range("H3").Select() | Out-Null
